$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(920, 1).Value = "Buying Opportunity"
$ws.Cells.Item(920, 2).Value = "support Zone"
$ws.Cells.Item(920, 3).Value = "long buildup"
$ws.Cells.Item(920, 4).Value = "Short buildup"
$ws.Cells.Item(920, 5).Value = "FII ENTERING"

$ws.Cells.Item(921, 1).Value = "AARVEEDEN"
$ws.Cells.Item(921, 2).Value = "ANGELONE"
$ws.Cells.Item(921, 3).Value = "DIXON"
$ws.Cells.Item(921, 6).Value = 27.37
$ws.Cells.Item(921, 7).Value = 2602
$ws.Cells.Item(921, 8).Value = 11242.85

$ws.Cells.Item(922, 1).Value = "ABSLNN50ET"
$ws.Cells.Item(922, 2).Value = "HINDMOTORS"
$ws.Cells.Item(922, 6).Value = 74.32
$ws.Cells.Item(922, 7).Value = 35.35

$ws.Cells.Item(923, 1).Value = "AKI"
$ws.Cells.Item(923, 2).Value = "MARUTI"
$ws.Cells.Item(923, 6).Value = 28.18
$ws.Cells.Item(923, 7).Value = 12845.2

$ws.Cells.Item(924, 1).Value = "ASTRAZEN"
$ws.Cells.Item(924, 2).Value = "MITTAL"
$ws.Cells.Item(924, 6).Value = 6458.55
$ws.Cells.Item(924, 7).Value = 2.23

$ws.Cells.Item(925, 1).Value = "AXISHCETF"
$ws.Cells.Item(925, 2).Value = "ORIENTALTL"
$ws.Cells.Item(925, 6).Value = 127.2
$ws.Cells.Item(925, 7).Value = 11.13

$ws.Cells.Item(926, 1).Value = "BEL"
$ws.Cells.Item(926, 6).Value = 309.6

$ws.Cells.Item(927, 1).Value = "BIKAJI"
$ws.Cells.Item(927, 6).Value = 725.55

$ws.Cells.Item(928, 1).Value = "COCHINSHIP"
$ws.Cells.Item(928, 6).Value = 2122.35

$ws.Cells.Item(929, 1).Value = "DCMSRIND"
$ws.Cells.Item(929, 6).Value = 236.64

$ws.Cells.Item(930, 1).Value = "DCW"
$ws.Cells.Item(930, 6).Value = 59.03

$ws.Cells.Item(931, 1).Value = "DCXINDIA"
$ws.Cells.Item(931, 6).Value = 359.85

$ws.Cells.Item(932, 1).Value = "DELHIVERY"
$ws.Cells.Item(932, 6).Value = 406.45

$ws.Cells.Item(933, 1).Value = "DHANI"
$ws.Cells.Item(933, 6).Value = 50.02

$ws.Cells.Item(934, 1).Value = "DTIL"
$ws.Cells.Item(934, 6).Value = 216.6

$ws.Cells.Item(935, 1).Value = "EMAMIPAP"
$ws.Cells.Item(935, 6).Value = 119.79

$ws.Cells.Item(936, 1).Value = "ESAFSFB"
$ws.Cells.Item(936, 6).Value = 53.54

$ws.Cells.Item(937, 1).Value = "GENUSPAPER"
$ws.Cells.Item(937, 6).Value = 21.77

$ws.Cells.Item(938, 1).Value = "GRSE"
$ws.Cells.Item(938, 6).Value = 1630

$ws.Cells.Item(939, 1).Value = "HAL"
$ws.Cells.Item(939, 6).Value = 5200.55

$ws.Cells.Item(940, 1).Value = "INOXGREEN"
$ws.Cells.Item(940, 6).Value = 141.51

$ws.Cells.Item(941, 1).Value = "IVP"
$ws.Cells.Item(941, 6).Value = 218.14

$ws.Cells.Item(942, 1).Value = "JSWINFRA"
$ws.Cells.Item(942, 6).Value = 300.4

$ws.Cells.Item(943, 1).Value = "KEC"
$ws.Cells.Item(943, 6).Value = 935.5

$ws.Cells.Item(944, 1).Value = "KINGFA"
$ws.Cells.Item(944, 6).Value = 2165.15

$ws.Cells.Item(945, 1).Value = "KMSUGAR"
$ws.Cells.Item(945, 6).Value = 48.07

$ws.Cells.Item(946, 1).Value = "KOHINOOR"
$ws.Cells.Item(946, 6).Value = 44.18

$ws.Cells.Item(947, 1).Value = "KRBL"
$ws.Cells.Item(947, 6).Value = 310.2

$ws.Cells.Item(948, 1).Value = "KSHITIJPOL"
$ws.Cells.Item(948, 6).Value = 7.7

$ws.Cells.Item(949, 1).Value = "KUANTUM"
$ws.Cells.Item(949, 6).Value = 159.91

$ws.Cells.Item(950, 1).Value = "MANINDS"
$ws.Cells.Item(950, 6).Value = 425.9

$ws.Cells.Item(951, 1).Value = "MCLEODRUSS"
$ws.Cells.Item(951, 6).Value = 29.3

$ws.Cells.Item(952, 1).Value = "MIDHANI"
$ws.Cells.Item(952, 6).Value = 457.55

$ws.Cells.Item(953, 1).Value = "MTARTECH"
$ws.Cells.Item(953, 6).Value = 1862.8

$ws.Cells.Item(954, 1).Value = "NAHARPOLY"
$ws.Cells.Item(954, 6).Value = 225.95

$ws.Cells.Item(955, 1).Value = "NIRAJ"
$ws.Cells.Item(955, 6).Value = 49.79

$ws.Cells.Item(956, 1).Value = "NUVOCO"
$ws.Cells.Item(956, 6).Value = 366.05

$ws.Cells.Item(957, 1).Value = "OAL"
$ws.Cells.Item(957, 6).Value = 378.8

$ws.Cells.Item(958, 1).Value = "ONEPOINT"
$ws.Cells.Item(958, 6).Value = 58.81

$ws.Cells.Item(959, 1).Value = "PARAS"
$ws.Cells.Item(959, 6).Value = 1156.9

$ws.Cells.Item(960, 1).Value = "PTCIL"
$ws.Cells.Item(960, 6).Value = 14729.6

$ws.Cells.Item(961, 1).Value = "RML"
$ws.Cells.Item(961, 6).Value = 895.8

$ws.Cells.Item(962, 1).Value = "14/06/2024"
